# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Septiembre de 2020 a las 15:52"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 6789877
$ws.Range("C4").Value = 1730
$ws.Range("D4").Value = 4069054
$ws.Range("E4").Value = 2520543
$ws.Range("G4").Value = 83
$ws.Range("H4").Value = 200280

# --- Row 5: India ---
$ws.Range("B5").Value = 5041681
$ws.Range("C5").Value = 23647
$ws.Range("D5").Value = 3960965
$ws.Range("E5").Value = 998430
$ws.Range("G5").Value = 195
$ws.Range("H5").Value = 82286

# --- Row 13: Argentina ---
$ws.Range("D13").Value = 448263
$ws.Range("E13").Value = 117165
$ws.Range("G13").Value = 58
$ws.Range("H13").Value = 11910

# --- Row 19: Arabia Saudita ---
$ws.Range("B19").Value = 327551
$ws.Range("C19").Value = 621
$ws.Range("D19").Value = 306004
$ws.Range("E19").Value = 17178
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 4369

# --- Row 59: Uzbekistan ---
$ws.Range("B59").Value = 49015
$ws.Range("C59").Value = 586
$ws.Range("D59").Value = 45422
$ws.Range("E59").Value = 3186
$ws.Range("G59").Value = 5
$ws.Range("H59").Value = 407

# --- Row 61: Suiza ---
$ws.Range("E61").Value = 6329
$ws.Range("G61").Value = 8
$ws.Range("H61").Value = 2036

# --- Rows 71/72: swap Serbia <-> Estado de Palestina (with updated data) ---
$ws.Range("A71").Value = "Estado de Palestina"
$ws.Range("B71").Value = 33006
$ws.Range("C71").Value = 756
$ws.Range("D71").Value = 22209
$ws.Range("E71").Value = 10554
$ws.Range("G71").Value = 14
$ws.Range("H71").Value = 243

$ws.Range("A72").Value = "Serbia"
$ws.Range("B72").Value = 32613
$ws.Range("C72").Value = 102
$ws.Range("D72").Value = 31411
$ws.Range("E72").Value = 466
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 736

# --- Rows 78/79: swap Bosnia y Herzegovina <-> Libia (with updated data) ---
$ws.Range("A78").Value = "Libia"
$ws.Range("B78").Value = 24936
$ws.Range("C78").Value = 792
$ws.Range("D78").Value = 13498
$ws.Range("E78").Value = 11044
$ws.Range("H78").Value = 394

$ws.Range("A79").Value = "Bosnia y Herzegovina"
$ws.Range("B79").Value = 24211
$ws.Range("C79").Value = 282
$ws.Range("D79").Value = 16990
$ws.Range("E79").Value = 6485
$ws.Range("G79").Value = 11
$ws.Range("H79").Value = 736

# --- Row 90: Zambia ---
$ws.Range("B90").Value = 13887
$ws.Range("C90").Value = 68
$ws.Range("D90").Value = 12869
$ws.Range("E90").Value = 692
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 326

# --- Row 101: Tayikistan ---
$ws.Range("B101").Value = 9171
$ws.Range("C101").Value = 42
$ws.Range("D101").Value = 7941
$ws.Range("E101").Value = 1157

# --- Row 136: Sri Lanka ---
$ws.Range("D136").Value = 3021
$ws.Range("E136").Value = 237

# --- Rows 204/205: swap Timor Oriental <-> Santa Lucia (data values are identical) ---
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Timor Oriental"

# --- Rows 214/215: swap Islas Malvinas <-> Montserrat ---
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
